# fix(backend): fix export column's order
#
# The "Colis ..." (parcels) block (rows 139-141) and the "Courriers ..."
# (mail) block (rows 142-144) were in the wrong order under the
# "Axe 3. Total des interactions du ..." section. Swap the two 3-row
# blocks so "Courriers" comes first, then "Colis" - values only, the
# row styles/formatting stay untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the current ("Colis" first, "Courriers" second) block values.
$colisEnregistres    = $ws.Cells.Item(139, 2).Value()
$colisRemis          = $ws.Cells.Item(140, 2).Value()
$colisReexpedies      = $ws.Cells.Item(141, 2).Value()
$courriersEnregistres = $ws.Cells.Item(142, 2).Value()
$courriersRemis       = $ws.Cells.Item(143, 2).Value()
$courriersReexpedies   = $ws.Cells.Item(144, 2).Value()

# Write them back swapped: "Courriers" block now first (139-141),
# "Colis" block now second (142-144).
$ws.Cells.Item(139, 2).Value = $courriersEnregistres
$ws.Cells.Item(140, 2).Value = $courriersRemis
$ws.Cells.Item(141, 2).Value = $courriersReexpedies

$ws.Cells.Item(142, 2).Value = $colisEnregistres
$ws.Cells.Item(143, 2).Value = $colisRemis
$ws.Cells.Item(144, 2).Value = $colisReexpedies

# Reflect the author's resulting selection/scroll position.
$null = $ws.Range("B142").Select()
